$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 289
$ws1.Range("F3").Value = 18
$ws1.Range("F4").Value = 7812
$ws1.Range("F5").Value = 5707
$ws1.Range("F6").Value = 474
$ws1.Range("F7").Value = 80
$ws1.Range("F10").Value = 264
$ws1.Range("F11").Value = 269
$ws1.Range("F12").Value = 58

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G3").Value = 80

# --- Sheet "全部类型" (all types, merged list) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 289
$ws4.Range("F3").Value = 18
$ws4.Range("F4").Value = 7812
$ws4.Range("F5").Value = 5707
$ws4.Range("F6").Value = 474
$ws4.Range("F7").Value = 80
$ws4.Range("F10").Value = 264
$ws4.Range("G12").Value = 80
$ws4.Range("F13").Value = 269
$ws4.Range("F14").Value = 58
